$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.449.91'
$ws.Range('E2').Value = '  +2.88%  '
$ws.Range('D3').Value = '3.381.69'
$ws.Range('E3').Value = '  +4.49%  '
$__style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '191.51'
$ws.Range('D5').Style = $__style
$ws.Range('E5').Value = '  +4.20%  '
$__style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '592.83'
$ws.Range('D6').Style = $__style
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.50%  '
$__style = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.134'
$ws.Range('D9').Style = $__style
$ws.Range('E9').Value = '  +2.96%  '
$ws.Range('E10').Value = '  +3.00%  '
$__style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.420'
$ws.Range('D11').Style = $__style
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('D12').Value = '3.973.90'
$ws.Range('E12').Value = '  +4.79%  '
$ws.Range('E13').Value = '  +1.26%  '
$__style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.68'
$ws.Range('D14').Style = $__style
$ws.Range('E14').Value = '  +3.80%  '
$ws.Range('D15').Value = '69.516.13'
$ws.Range('E15').Value = '  +2.92%  '
$ws.Range('E16').Value = '  +1.96%  '
$ws.Range('D17').Value = '3.387.56'
$ws.Range('E17').Value = '  +5.15%  '
$__style = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '449.17'
$ws.Range('D18').Style = $__style
$ws.Range('E18').Value = '  +13.80%  '
$__style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.85'
$ws.Range('D19').Style = $__style
$ws.Range('E19').Value = '  +1.61%  '
$__style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.82'
$ws.Range('D20').Style = $__style
$ws.Range('E20').Value = '  +2.57%  '
$__style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.82'
$ws.Range('D21').Style = $__style
$ws.Range('E21').Value = '  +3.49%  '
$__style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.13'
$ws.Range('D22').Style = $__style
$ws.Range('E22').Value = '  +5.35%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').Value = '3.523.33'
$ws.Range('E24').Value = '  +4.44%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$__style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.522'
$ws.Range('D25').Style = $__style
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$__style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000123'
$ws.Range('D26').Style = $__style
$ws.Range('E26').Value = '  +4.05%  '
$ws.Range('E27').Value = '  +1.70%  '
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +1.92%  '
$__style = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.44'
$ws.Range('D31').Style = $__style
$ws.Range('E31').Value = '  +3.70%  '
$__style = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.66'
$ws.Range('D32').Style = $__style
$ws.Range('E32').Value = '  +1.84%  '
$ws.Range('E33').Value = '  +3.16%  '
$__style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.99'
$ws.Range('D34').Style = $__style
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  +5.42%  '
$__style = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '165.37'
$ws.Range('D37').Style = $__style
$ws.Range('E37').Value = '  +2.96%  '
$__style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.95'
$ws.Range('D38').Style = $__style
$ws.Range('E38').Value = '  +3.44%  '
$__style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.49'
$ws.Range('D39').Style = $__style
$ws.Range('E39').Value = '  +4.11%  '
$__style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.818'
$ws.Range('D40').Style = $__style
$ws.Range('E40').Value = '  +1.90%  '
$__style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.59'
$ws.Range('D41').Style = $__style
$ws.Range('E41').Value = '  +1.01%  '
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('D43').Value = '2.756.63'
$ws.Range('E43').Value = '  +5.49%  '
$__style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.54'
$ws.Range('D44').Style = $__style
$ws.Range('E44').Value = '  +2.85%  '
$__style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.59'
$ws.Range('D45').Style = $__style
$ws.Range('E45').Value = '  +3.48%  '
$ws.Range('E46').Value = '  +0.75%  '
$__style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.82'
$ws.Range('D47').Style = $__style
$ws.Range('E47').Value = '  +0.78%  '
$__style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '340.97'
$ws.Range('D48').Style = $__style
$ws.Range('E48').Value = '  +2.13%  '
$ws.Range('E49').Value = '  +2.13%  '
$__style = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.09'
$ws.Range('D50').Style = $__style
$ws.Range('E50').Value = '  +8.14%  '
$ws.Range('E51').Value = '  +6.08%  '
